$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 ("2025") metrics with revised figures
$ws.Range("C8").Value = 1225
$ws.Range("D8").Value = 199
$ws.Range("E8").Value = 1026
$ws.Range("F8").Value = 8.162428219852339
$ws.Range("G8").Value = 83.75510204081633
$ws.Range("H8").Value = 16.24489795918367
